$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.813.93"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "2.349.59"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'0.671"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").Value = "'237.26"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("D7").Value = "'73.00"
$ws.Range("E7").Value = "  +10.70%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +19.12%  "
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("D11").Value = "'28.98"
$ws.Range("E11").Value = "  +8.13%  "
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").Value = "2.697.75"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "'16.71"
$ws.Range("E14").Value = "  +8.73%  "
$ws.Range("E15").Value = "  +6.37%  "
$ws.Range("D16").Value = "'0.896"
$ws.Range("E16").Value = "  +5.76%  "
$ws.Range("D17").Value = "2.338.87"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "43.801.20"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  +3.50%  "
$ws.Range("D20").Value = "'77.72"
$ws.Range("E20").Value = "  +4.85%  "
$ws.Range("D21").Value = "'6.41"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").Value = "'253.75"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("D26").Value = "'10.52"
$ws.Range("E26").Value = "  +6.13%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'22.37"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "'172.45"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  +6.27%  "
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +4.83%  "
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("D34").Value = "'0.0716"
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("D35").Value = "'5.21"
$ws.Range("E35").Value = "  +5.11%  "
$ws.Range("D36").Value = "'3.99"
$ws.Range("E36").Value = "  +10.25%  "
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("D38").Value = "'6.38"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").Value = "'0.0266"
$ws.Range("E39").Value = "  +6.01%  "
$ws.Range("D40").Value = "'19.57"
$ws.Range("E40").Value = "  +7.90%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'8.81"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").Value = "'1.25"
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").Value = "'4.44"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").Value = "'97.98"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  +11.01%  "
$ws.Range("D49").Value = "'2.33"
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("D50").Value = "1.434.12"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("E51").Value = "  +1.52%  "
